$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Archive
#
# The handoff/handback bookkeeping rows for the three files
#   48cbe1ab-a6ea-4d2a-af77-62167ea810bf.md
#   1188accf-185d-45a8-9d6a-6a871890c400.md
#   1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.md
# are re-ordered (1188accf, 1d5cb4a2, 48cbe1ab) across the Overview, zh-cn
# and de-de sheets, and the "48cbe1ab" entry picks up a new "In Translation"
# status/date now that 1188accf/1d5cb4a2 have moved on to translation.
# ---------------------------------------------------------------------------

# ---- Overview sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A4").Value = "1188accf-185d-45a8-9d6a-6a871890c400.md"
$ws.Range("B4").Value = "e2e\1188accf-185d-45a8-9d6a-6a871890c400.md"
$ws.Range("G4").Value = "2016-08-24 10:18:49"

$ws.Range("A5").Value = "1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.md"
$ws.Range("B5").Value = "e2e\1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.md"
$ws.Range("E5").Value = "In Translation"
$ws.Range("F5").Value = "In Translation"

$ws.Range("A6").Value = "48cbe1ab-a6ea-4d2a-af77-62167ea810bf.md"
$ws.Range("B6").Value = "e2e\48cbe1ab-a6ea-4d2a-af77-62167ea810bf.md"
$ws.Range("E6").Value = "In Translation"
$ws.Range("F6").Value = "In Translation"
$ws.Range("G6").Value = "2016-08-24 10:13:38"

# ---- zh-cn sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A4").Value = "1188accf-185d-45a8-9d6a-6a871890c400.md"
$ws.Range("E4").Value = ""
$ws.Range("G4").Value = "1188accf-185d-45a8-9d6a-6a871890c400.f73616abc37196bca00c69555d1b91353c4079d9.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-24 10:18:44"

$ws.Range("A5").Value = "1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.f6b758931d74057829e64593d9ce7c4dbc6e0130.zh-cn.xlf"

$ws.Range("A6").Value = "48cbe1ab-a6ea-4d2a-af77-62167ea810bf.md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("E6").Value = "ht"
$ws.Range("G6").Value = "48cbe1ab-a6ea-4d2a-af77-62167ea810bf.723f22c7418630b44ac2c076a89dafbd77a8c454.zh-cn.xlf"
$ws.Range("H6").Value = "2016-08-24 10:13:33"

# ---- de-de sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A4").Value = "1188accf-185d-45a8-9d6a-6a871890c400.md"
$ws.Range("E4").Value = ""
$ws.Range("G4").Value = "1188accf-185d-45a8-9d6a-6a871890c400.f73616abc37196bca00c69555d1b91353c4079d9.de-de.xlf"
$ws.Range("H4").Value = "2016-08-24 10:18:49"

$ws.Range("A5").Value = "1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("G5").Value = "1d5cb4a2-4a5c-4cb8-96c3-d4b0ac96a9e2.f6b758931d74057829e64593d9ce7c4dbc6e0130.de-de.xlf"

$ws.Range("A6").Value = "48cbe1ab-a6ea-4d2a-af77-62167ea810bf.md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("E6").Value = "ht"
$ws.Range("G6").Value = "48cbe1ab-a6ea-4d2a-af77-62167ea810bf.723f22c7418630b44ac2c076a89dafbd77a8c454.de-de.xlf"
$ws.Range("H6").Value = "2016-08-24 10:13:38"
